# Natmi following Dr Hou advice
# Update Cthrc1-Ror2 LR-pair sheet: revise row 2 and add rows 3-7 for the
# full Sending-cluster x Target-cluster combination grid (FAPs/sCs senders,
# FAPs/sCs/ECs targets).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cthrc1"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.449420333333333
$ws.Range("H2").Value = 16.348261
$ws.Range("I2").Value = 0.8985142489564721
$ws.Range("J2").Value = 0.8985142489564723
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.028814
$ws.Range("N2").Value = 0.086442
$ws.Range("O2").Value = 0.009404897244300481
$ws.Range("P2").Value = 0.009404897244300481
$ws.Range("Q2").Value = 0.1570195974846667
$ws.Range("R2").Value = 1.413176377362
$ws.Range("S2").Value = 0.008450434183975442
$ws.Range("T2").Value = 0.008450434183975442

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cthrc1"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.449420333333333
$ws.Range("H3").Value = 16.348261
$ws.Range("I3").Value = 0.8985142489564721
$ws.Range("J3").Value = 0.8985142489564723
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.294804333333333
$ws.Range("N3").Value = 6.884412999999999
$ws.Range("O3").Value = 0.7490247432073112
$ws.Range("P3").Value = 0.7490247432073114
$ws.Range("Q3").Value = 12.50535339508811
$ws.Range("R3").Value = 112.548180555793
$ws.Range("S3").Value = 0.6730094045927316
$ws.Range("T3").Value = 0.6730094045927318

# Row 4: FAPs -> sCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cthrc1"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.449420333333333
$ws.Range("H4").Value = 16.348261
$ws.Range("I4").Value = 0.8985142489564721
$ws.Range("J4").Value = 0.8985142489564723
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7401046666666667
$ws.Range("N4").Value = 2.220314
$ws.Range("O4").Value = 0.2415703595483883
$ws.Range("P4").Value = 0.2415703595483883
$ws.Range("Q4").Value = 4.033141419328222
$ws.Range("R4").Value = 36.298272773954
$ws.Range("S4").Value = 0.217054410179765
$ws.Range("T4").Value = 0.2170544101797651

# Row 5: sCs -> ECs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Cthrc1"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.6155033333333333
$ws.Range("H5").Value = 1.84651
$ws.Range("I5").Value = 0.1014857510435278
$ws.Range("J5").Value = 0.1014857510435278
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.028814
$ws.Range("N5").Value = 0.086442
$ws.Range("O5").Value = 0.009404897244300481
$ws.Range("P5").Value = 0.009404897244300481
$ws.Range("Q5").Value = 0.01773511304666667
$ws.Range("R5").Value = 0.15961601742
$ws.Range("S5").Value = 0.0009544630603250395
$ws.Range("T5").Value = 0.0009544630603250397

# Row 6: sCs -> FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Cthrc1"
$ws.Range("C6").Value = "Ror2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.6155033333333333
$ws.Range("H6").Value = 1.84651
$ws.Range("I6").Value = 0.1014857510435278
$ws.Range("J6").Value = 0.1014857510435278
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.294804333333333
$ws.Range("N6").Value = 6.884412999999999
$ws.Range("O6").Value = 0.7490247432073112
$ws.Range("P6").Value = 0.7490247432073114
$ws.Range("Q6").Value = 1.412459716514444
$ws.Range("R6").Value = 12.71213744863
$ws.Range("S6").Value = 0.07601533861457954
$ws.Range("T6").Value = 0.07601533861457957

# Row 7: sCs -> sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Cthrc1"
$ws.Range("C7").Value = "Ror2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.6155033333333333
$ws.Range("H7").Value = 1.84651
$ws.Range("I7").Value = 0.1014857510435278
$ws.Range("J7").Value = 0.1014857510435278
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7401046666666667
$ws.Range("N7").Value = 2.220314
$ws.Range("O7").Value = 0.2415703595483883
$ws.Range("P7").Value = 0.2415703595483883
$ws.Range("Q7").Value = 0.4555368893488889
$ws.Range("R7").Value = 4.09983200414
$ws.Range("S7").Value = 0.02451594936862323
$ws.Range("T7").Value = 0.02451594936862324
